# Generate Report for Handoff
#
# A new localization handoff pass was generated for 4 files
# (285bbb65-..., 2db416ef-..., b486bf50-..., f0d158b5-...) which were
# previously at low priority and pending handoff. Update their Priority
# to "ht" and refresh the Latest Handoff Datetime on both locale sheets,
# and mirror the de-de handoff date into the Overview sheet's
# "Latest HO Xliff Generate Date" column (these cells track the de-de
# handoff date for these rows).

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4:E7").Value = "ht"
$wsZhCn.Range("H4:H7").Value = "2016-08-18 00:29:31"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4:E7").Value = "ht"
$wsDeDe.Range("H4:H7").Value = "2016-08-18 00:29:36"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4:G7").Value = "2016-08-18 00:29:36"
